$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Cell content changes -------------------------------------------------

# Row 2: "NewTurn" -> "Intel", "RaceTurn" -> "Orders"
$ws.Range("J2").Value = "Intel"
$ws.Range("L2").Value = "Orders"

# New mirrored cells in column I for rows 4-6 (match column B/H on same row)
$ws.Range("I4").Value = $ws.Range("H4").Value
$ws.Range("I5").Value = $ws.Range("H5").Value
$ws.Range("I6").Value = $ws.Range("H6").Value

# Row 9: "AllStars (position only)" -> "AllStars"
$ws.Range("I9").Value = "AllStars"

# Row 39: "RaceData" -> "RaceData (not a RaceData object)"
$ws.Range("D39").Value = "RaceData (not a RaceData object)"

# --- View / selection changes ---------------------------------------------
# The author had scrolled the window down (topLeftCell A16) and left the
# active selection on D40.
$ws.Range("D40").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
